$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "MEC-2B-Máquinas Térmicas e de Fluxo" from Tuesday (C4, C6) to Friday (F7, F8)
$ws.Range("C4").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("F7").Value = "MEC-2B-Máquinas Térmicas e de Fluxo"
$ws.Range("F8").Value = "MEC-2B-Máquinas Térmicas e de Fluxo"
